# Auto-generated edits for Siren_Profits workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 62500716
$ws.Range("I9").Value = 100000300
$ws.Range("K9").Value = 100000300
$ws.Range("M9").Value = -100000131

$ws.Range("H18").Value = 7290.591
$ws.Range("I18").Value = 7519.7
$ws.Range("J18").Value = 4999.5
$ws.Range("K18").Value = 7519.7
$ws.Range("L18").Value = 4999.5
$ws.Range("M18").Value = -7235.7
$ws.Range("N18").Value = -5567.5

$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H62").Value = 4442.467
$ws.Range("I62").Value = 4196.1816
$ws.Range("K62").Value = 4196.1816
$ws.Range("M62").Value = -3572.1816

$ws.Range("H65").Value = 4442.467
$ws.Range("I65").Value = 4196.1816
$ws.Range("K65").Value = 20980.908
$ws.Range("M65").Value = -17860.908

$ws.Range("H101").Value = 10992484
$ws.Range("I101").Value = 14289005
$ws.Range("K101").Value = 42867015
$ws.Range("M101").Value = -42865393

$ws.Range("H132").Value = 14943.211
$ws.Range("I132").Value = 17751.643
$ws.Range("J132").Value = 7079.6
$ws.Range("K132").Value = 53254.929
$ws.Range("L132").Value = 21238.8
$ws.Range("M132").Value = -50724.929
$ws.Range("N132").Value = -26298.8

$ws.Range("H138").Value = 3696.8518
$ws.Range("I138").Value = 3303.25
$ws.Range("J138").Value = 3862.5789
$ws.Range("K138").Value = 9909.75
$ws.Range("L138").Value = 11587.7367
$ws.Range("M138").Value = -4769.75
$ws.Range("N138").Value = -21867.7367

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5796.704
$ws.Range("I32").Value = 6337.6553
$ws.Range("K32").Value = 6337.6553
$ws.Range("M32").Value = -6050.6553

$ws.Range("H45").Value = 8530.875
$ws.Range("I45").Value = 10983
$ws.Range("K45").Value = 10983
$ws.Range("M45").Value = -10606

$ws.Range("H61").Value = 3906.8813
$ws.Range("I61").Value = 3900.1072
$ws.Range("J61").Value = 4033.3333
$ws.Range("K61").Value = 3900.1072
$ws.Range("L61").Value = 4033.3333
$ws.Range("M61").Value = -3688.1072
$ws.Range("N61").Value = -4457.3333

$ws.Range("H122").Value = 1226373
$ws.Range("I122").Value = 3927.65
$ws.Range("J122").Value = 2754429.8
$ws.Range("K122").Value = 11782.95
$ws.Range("L122").Value = 8263289.399999999
$ws.Range("M122").Value = -9332.950000000001
$ws.Range("N122").Value = -8268189.399999999

$ws.Range("H136").Value = 3906.8813
$ws.Range("I136").Value = 3900.1072
$ws.Range("J136").Value = 4033.3333
$ws.Range("K136").Value = 11700.3216
$ws.Range("L136").Value = 12099.9999
$ws.Range("M136").Value = -9150.321599999999
$ws.Range("N136").Value = -17199.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3799.7144
$ws.Range("I86").Value = 4265.4585
$ws.Range("J86").Value = 2783.5454
$ws.Range("K86").Value = 4265.4585
$ws.Range("L86").Value = 2783.5454
$ws.Range("M86").Value = -3142.4585
$ws.Range("N86").Value = -5029.5454

$ws.Range("H89").Value = 3799.7144
$ws.Range("I89").Value = 4265.4585
$ws.Range("J89").Value = 2783.5454
$ws.Range("K89").Value = 21327.2925
$ws.Range("L89").Value = 13917.727
$ws.Range("M89").Value = -15711.2925
$ws.Range("N89").Value = -25149.727

$ws.Range("H134").Value = 17358.188
$ws.Range("I134").Value = 17358.188
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 52074.564
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -49539.564
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4664.077
$ws.Range("I31").Value = 4261.0835
$ws.Range("K31").Value = 4261.0835
$ws.Range("M31").Value = -3966.0835

$ws.Range("H34").Value = 4664.077
$ws.Range("I34").Value = 4261.0835
$ws.Range("K34").Value = 4261.0835
$ws.Range("M34").Value = -4059.0835

$ws.Range("H43").Value = 131885.67
$ws.Range("J43").Value = 131885.67
$ws.Range("L43").Value = 131885.67
$ws.Range("N43").Value = -132253.67

$ws.Range("H58").Value = 2067.8965
$ws.Range("I58").Value = 962.2941
$ws.Range("K58").Value = 962.2941
$ws.Range("M58").Value = -759.2941

$ws.Range("H101").Value = 131885.67
$ws.Range("J101").Value = 131885.67
$ws.Range("L101").Value = 131885.67
$ws.Range("N101").Value = -138375.67

$ws.Range("H132").Value = 70491.164
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H134").Value = 4124.385
$ws.Range("I134").Value = 3801.7144
$ws.Range("K134").Value = 11405.1432
$ws.Range("M134").Value = -8870.143199999999

$ws.Range("H136").Value = 2067.8965
$ws.Range("I136").Value = 962.2941
$ws.Range("K136").Value = 2886.8823
$ws.Range("M136").Value = -336.8822999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 11867.083
$ws.Range("I68").Value = 2249.3333
$ws.Range("K68").Value = 6747.999899999999
$ws.Range("M68").Value = -5936.999899999999

$ws.Range("H71").Value = 11867.083
$ws.Range("I71").Value = 2249.3333
$ws.Range("K71").Value = 20243.9997
$ws.Range("M71").Value = -16187.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 16094.091
$ws.Range("I102").Value = 17003.5
$ws.Range("K102").Value = 17003.5
$ws.Range("M102").Value = -15381.5

$ws.Range("H107").Value = 424.85715
$ws.Range("I107").Value = 514.36365
$ws.Range("J107").Value = 96.666664
$ws.Range("K107").Value = 514.36365
$ws.Range("L107").Value = 96.666664
$ws.Range("M107").Value = 1405.63635
$ws.Range("N107").Value = -3936.666664

$ws.Range("H122").Value = 8529.259
$ws.Range("I122").Value = 4879.6313
$ws.Range("K122").Value = 14638.8939
$ws.Range("M122").Value = -12188.8939

$ws.Range("H138").Value = 149999.5
$ws.Range("J138").Value = 149999.5
$ws.Range("L138").Value = 149999.5
$ws.Range("N138").Value = -160279.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17229.572
$ws.Range("I7").Value = 30425.8
$ws.Range("K7").Value = 30425.8
$ws.Range("M7").Value = -30313.8

$ws.Range("H16").Value = 3315.9473
$ws.Range("I16").Value = 3315.9473
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3315.9473
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3145.9473
$ws.Range("N16").ClearContents()

$ws.Range("H46").Value = 3166.6316
$ws.Range("I46").Value = 847.25
$ws.Range("K46").Value = 847.25
$ws.Range("M46").Value = -659.25

$ws.Range("H61").Value = 4506.769
$ws.Range("I61").Value = 2814.3684
$ws.Range("K61").Value = 2814.3684
$ws.Range("M61").Value = -2612.3684

$ws.Range("H100").Value = 2295.5454
$ws.Range("I100").Value = 2199
$ws.Range("K100").Value = 2199
$ws.Range("M100").Value = -1658

$ws.Range("H113").Value = 4506.769
$ws.Range("I113").Value = 2814.3684
$ws.Range("K113").Value = 2814.3684
$ws.Range("M113").Value = -644.3683999999998

$ws.Range("H122").Value = 4850.968
$ws.Range("I122").Value = 4620.7144
$ws.Range("K122").Value = 13862.1432
$ws.Range("M122").Value = -11412.1432

$ws.Range("H126").Value = 17229.572
$ws.Range("I126").Value = 30425.8
$ws.Range("K126").Value = 91277.39999999999
$ws.Range("M126").Value = -88807.39999999999

$ws.Range("H132").Value = 600943.0600000001
$ws.Range("I132").Value = 1865112.6
$ws.Range("J132").Value = 6039.706
$ws.Range("K132").Value = 5595337.800000001
$ws.Range("L132").Value = 18119.118
$ws.Range("M132").Value = -5592807.800000001
$ws.Range("N132").Value = -23179.118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 24751876
$ws.Range("I5").Value = 29200000
$ws.Range("K5").Value = 29200000
$ws.Range("M5").Value = -29199888

$ws.Range("H100").Value = 19723.834
$ws.Range("I100").Value = 15851.704
$ws.Range("J100").Value = 26693.666
$ws.Range("K100").Value = 31703.408
$ws.Range("L100").Value = 53387.332
$ws.Range("M100").Value = -31162.408
$ws.Range("N100").Value = -54469.332

$ws.Range("H122").Value = 27797.592
$ws.Range("I122").Value = 4567.7144
$ws.Range("J122").Value = 68449.875
$ws.Range("K122").Value = 13703.1432
$ws.Range("L122").Value = 205349.625
$ws.Range("M122").Value = -11253.1432
$ws.Range("N122").Value = -210249.625

$ws.Range("H126").Value = 25113.666
$ws.Range("I126").Value = 28717.133
$ws.Range("K126").Value = 86151.399
$ws.Range("M126").Value = -83681.399

$ws.Range("H132").Value = 5710.108
$ws.Range("I132").Value = 5631.7017
$ws.Range("K132").Value = 16895.1051
$ws.Range("M132").Value = -14365.1051

$ws.Range("H135").Value = 10000.5
$ws.Range("J135").Value = 10000.5
$ws.Range("L135").Value = 10000.5
$ws.Range("N135").Value = -20140.5

$ws.Range("H136").Value = 276151.84
$ws.Range("I136").Value = 321508.4
$ws.Range("J136").Value = 4012.375
$ws.Range("K136").Value = 964525.2000000001
$ws.Range("L136").Value = 12037.125
$ws.Range("M136").Value = -961975.2000000001
$ws.Range("N136").Value = -17137.125
